# Merge completion: mark additional rows as "merged" (and flag a subset with "*")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "merged"
$ws.Range("D3").Value = "merged"
$ws.Range("D4").Value = "merged"
$ws.Range("D12").Value = "merged"
$ws.Range("D18").Value = "merged"
$ws.Range("D22").Value = "merged"
$ws.Range("D24").Value = "merged"
$ws.Range("D30").Value = "merged"
$ws.Range("D32").Value = "merged"
$ws.Range("D34").Value = "merged"
$ws.Range("E34").Value = "*"
$ws.Range("D35").Value = "merged"
$ws.Range("E35").Value = "*"
$ws.Range("D38").Value = "merged"
$ws.Range("E38").Value = "*"
$ws.Range("D40").Value = "merged"
$ws.Range("D42").Value = "merged"
$ws.Range("E42").Value = "*"
$ws.Range("D43").Value = "merged"
$ws.Range("D44").Value = "merged"
$ws.Range("D46").Value = "merged"
$ws.Range("D59").Value = "merged"
$ws.Range("E59").Value = "*"
$ws.Range("D61").Value = "merged"
$ws.Range("D67").Value = "merged"
$ws.Range("D68").Value = "merged"
$ws.Range("D69").Value = "merged"
$ws.Range("E76").Value = "*"
$ws.Range("E77").Value = "*"
$ws.Range("D78").Value = "merged"
$ws.Range("E78").Value = "*"
$ws.Range("D79").Value = "merged"
$ws.Range("D87").Value = "merged"
$ws.Range("E87").Value = "*"
$ws.Range("D88").Value = "merged"
$ws.Range("E88").Value = "*"
$ws.Range("D92").Value = "merged"
$ws.Range("D94").Value = "merged"
$ws.Range("D96").Value = "merged"
$ws.Range("D97").Value = "merged"
$ws.Range("D98").Value = "merged"
$ws.Range("D100").Value = "merged"
$ws.Range("E100").Value = "*"
$ws.Range("D101").Value = "merged"
$ws.Range("E101").Value = "*"
$ws.Range("D102").Value = "merged"
$ws.Range("D104").Value = "merged"
$ws.Range("D108").Value = "merged"
$ws.Range("D112").Value = "merged"
$ws.Range("D114").Value = "merged"
$ws.Range("D122").Value = "merged"
$ws.Range("D123").Value = "merged"
$ws.Range("D124").Value = "merged"
$ws.Range("D135").Value = "merged"
$ws.Range("D136").Value = "merged"
$ws.Range("D137").Value = "merged"
$ws.Range("D159").Value = "merged"
$ws.Range("D161").Value = "merged"
$ws.Range("F5").Value = "* denotes file requiring later attention"

# Re-color the "yes" conditional format to match the "unique" (yellow) style
$yesUniqueFcs = $ws.Range("B2:C165").FormatConditions
$yesCondition = $yesUniqueFcs.Item(1)
$yesCondition.Font.Color = 26012
$yesCondition.Interior.Color = 10284031

# Move the visible selection to where the merge work left off
$ws.Range("D162").Select()

